# Insert a new data row before the current row 131 (shifts existing rows
# 131-132 down to 132-133) and populate it with the new weekly price entry
# for Femacal de La Calera - Albahaca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(131).Insert()

$ws.Cells.Item(131, 1).Value2 = 3
$ws.Cells.Item(131, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(131, 3).Value2 = "Coquimbo"
$ws.Cells.Item(131, 4).Value2 = 44595
$ws.Cells.Item(131, 5).Value2 = 5
$ws.Cells.Item(131, 6).Value2 = 100112052
$ws.Cells.Item(131, 7).Value2 = "Albahaca"
$ws.Cells.Item(131, 8).Value2 = "Sin especificar"
$ws.Cells.Item(131, 9).Value2 = "Primera"
$ws.Cells.Item(131, 10).Value2 = 150
$ws.Cells.Item(131, 11).Value2 = 4000
$ws.Cells.Item(131, 12).Value2 = 4500
$ws.Cells.Item(131, 13).Value2 = 4233
$ws.Cells.Item(131, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(131, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(131, 16).Value2 = 706
$ws.Cells.Item(131, 17).Value2 = 6
$ws.Cells.Item(131, 18).Value2 = "Hortaliza"

# Make sure the style used for the new date cell matches the rest of the
# date column (numFmtId 165 / style index 2), same as D130/D132.
$ws.Cells.Item(131, 4).NumberFormat = $ws.Cells.Item(130, 4).NumberFormat
